$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has 14 data rows (row 2..15) that all share the exact
# same repositoryUrl / repositoryName / repositoryAuthor / startingDate and
# the same OSE/BCE/PDE/SV/OS/SD/RS/TFS/UI/TC text flags (all stored as plain
# text, e.g. "0"/"1", not numbers). Row 15 is used as a formatting/typing
# template: copying it down preserves the text typing (and cell styling) of
# every column instead of letting Excel reinterpret strings like
# "12/03/2018" or "0"/"1" as dates/numbers.
$templateRow = 15

# New rows 16-21 (ids 15-20), with the OS (J) column flag updated for some
# rows as introduced by the web-service JSON change.
$newRows = @(
    @{ Row = 16; Id = 15; J = "0" },
    @{ Row = 17; Id = 16; J = "1" },
    @{ Row = 18; Id = 17; J = "1" },
    @{ Row = 19; Id = 18; J = "1" },
    @{ Row = 20; Id = 19; J = "1" },
    @{ Row = 21; Id = 20; J = "0" }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $srcRange = $ws.Range($ws.Cells.Item($templateRow, 1), $ws.Cells.Item($templateRow, 15))
    $dstRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 15))
    $srcRange.Copy($dstRange)

    # Fix up the id (column A) for this row.
    $ws.Cells.Item($r, 1).Value = $rowData.Id

    # Fix up the OS (column J) flag when it differs from the template row,
    # copying from column I (which already holds a text "1") so the cell
    # keeps its text typing instead of becoming a number.
    if ($rowData.J -eq "1") {
        $ws.Cells.Item($templateRow, 9).Copy($ws.Cells.Item($r, 10))
    }
}
